$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / URL text updates (A1, A2) ---
$ws.Range("A1").Value = "Municipal Recycling Program Comparison: Data Dictionary"
$ws.Range("A2").Value = "URL: https://github.com/ewhinihan/municipal_recycling"

# --- Row 5 (city_Program) / Row 6 (year) description text updates ---
$ws.Range("F6").Value = "Year that recycle data was taken from."
$ws.Range("F5").Value = "Name of city/municipality for relevant recycle data."

# --- Reorder variable rows 11-17: move contaminents/tonnage up to follow the
#     tonnage-type rows, and push total_Cost/cost_Per_Ton/recycle_Rate/population/area down ---
$newRows = @{
    11 = @("Contaminated Materials Tonnage","contaminents","Integer","Whole number","Integers greater than 0","Total tonnage of contaminated material (from city data); usually means hazardous waste.","y","n")
    12 = @("Total number of tons recycled","tonnage","Integer","Whole number","Integers greater than 0","Amount recycled in tons for the given year in each city (cumulative residential totals).","y","n")
    13 = @("Total cost per ton","total_Cost","Integer","Whole number","Numbers greater than 0 (USD)","Total cost to recycle all materials in each city for given year.","n","y")
    14 = @("Cost per ton of material recycled","cost_Per_Ton","Integer","Decimal","Numbers greater than 0 (USD)","Amount of money it costs to recycle one ton of recyclable material.","n","y")
    15 = @("Recycle rate","recycle_Rate","Integer","Fraction","Numbers between 0-1 (percentage)","Total percentage of residential (includes single-family and multi-family) recycled materials that were recovered or diverted from a landfill. The numbers for Seattle were pulled direclty from report; Portland and Los Angeles had to be calculated by hand from available data.","y","n")
    16 = @("Population of City","population","Integer","Whole number","Integers greater than 0","Population of city measured in individual residents (US Census data).","y","n")
    17 = @("Total area of city","area","Integer","Number","Numbers greater than 0 (miles squared)","Area of each city measured in square miles (US Census data).","y","n")
}

foreach ($r in $newRows.Keys) {
    $rowVals = $newRows[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $rowVals[$i]
    }
}

# --- Column width tweaks ---
$ws.Columns.Item(2).ColumnWidth = 18.109375
$ws.Columns.Item(3).ColumnWidth = 16.44140625

# --- Row height tweaks ---
$ws.Rows.Item(11).RowHeight = 31.2
$ws.Rows.Item(14).RowHeight = 37.2
$ws.Rows.Item(15).RowHeight = 93

# --- Selection state ---
$ws.Range("F7").Select()

# --- Window view tweaks ---
$excel.ActiveWindow.WindowState = -4143
$excel.Width = 29040
$excel.Height = 15840
